$d = $word.ActiveDocument

# --- Step 1: locate the end of the existing last paragraph's text -------
# ("...to be easily saved in a file.") and split it into a new paragraph
# that will hold the new sentences, preserving the ListParagraph /
# ilvl=1 / numId=13 formatting (InsertParagraphAfter clones pPr).
$find = $d.Content
$find.Find.Execute("to be easily saved in a file.", $true, $false, $false, `
                    $false, $false, $true, 1, $false, "", 0)
$splitPos = $find.End
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# --- Step 2: type the new sentences into the freshly created paragraph --
# Each InsertAfter call below lands in its own run, mirroring the
# run-boundaries of the authored text (word/phrase groupings).
$pos = $splitPos + 1

$t = "SerialImport"
$r = $d.Range($pos, $pos)
$r.InsertAfter($t)
$pos = $pos + $t.Length

$t = " and "
$r = $d.Range($pos, $pos)
$r.InsertAfter($t)
$pos = $pos + $t.Length

$t = "SerialExport"
$r = $d.Range($pos, $pos)
$r.InsertAfter($t)
$pos = $pos + $t.Length

$t = " were put into separate classes "
$r = $d.Range($pos, $pos)
$r.InsertAfter($t)
$pos = $pos + $t.Length

$t = "so as to"
$r = $d.Range($pos, $pos)
$r.InsertAfter($t)
$pos = $pos + $t.Length

$t = " not give "
$r = $d.Range($pos, $pos)
$r.InsertAfter($t)
$pos = $pos + $t.Length

$t = "DimensionalSpace"
$r = $d.Range($pos, $pos)
$r.InsertAfter($t)
$pos = $pos + $t.Length

$t = " too many responsibilities. These classes also handle the exportation and importation of Point and Cell objects as well."
$r = $d.Range($pos, $pos)
$r.InsertAfter($t)
$pos = $pos + $t.Length

# $pos now sits right after "...Point and Cell objects as well." - this is
# where the _GoBack bookmark belongs (before the trailing space run).
$bmTarget = $pos

# --- Step 3: append a trailing space so the paragraph ends in "well. " --
$spaceRange = $d.Range($bmTarget, $bmTarget)
$spaceRange.InsertAfter(" ")

# --- Step 4: pad further out so the bookmark target position is not the
# very last position in the story (this engine mis-places zero-width
# bookmarks added too close to the end of the document), then trim the
# padding back out once the bookmark is safely anchored.
$padPos = $d.Paragraphs.Last.Range.End - 1
$padRange = $d.Range($padPos, $padPos)
$padRange.InsertAfter("PADDINGPADDINGPADDINGPADDING")

# --- Step 5: move the _GoBack bookmark to sit right after the new
# sentences and before the trailing space run.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Range($bmTarget, $bmTarget)
$d.Bookmarks.Add("_GoBack", $bmRange)
$bm = $d.Bookmarks("_GoBack")

# --- Step 6: remove the padding text now that the bookmark is placed ----
$lastPara = $d.Paragraphs.Last
$padDelete = $d.Range($bm.End + 1, $lastPara.Range.End - 1)
$padDelete.Delete()

Write-Output "Done"
